$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.198741555213928
$ws.Range("B1").Value = 2.012449502944946
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.03538179397583
$ws.Range("E1").Value = 1.208709955215454
